$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (also updates the ExternalData_1 defined name reference)
$ws.Name = "map"

# The Power Query refresh re-wrote the table body without the old
# "applyNumberFormat" style (cellXfs index 1) and dropped any cell that
# ended up completely blank (no value, no explicit formatting).
$ws.Range("B2:E110").ClearFormats() | Out-Null

$emptyCells = "C36,D36,E36,D46,D47,D61,D63,E63,D66,E66,D67,C75,D75,E75,D79,C83,D83,E83,C86,D86,E86,D88,D93,D94"
foreach ($addr in $emptyCells.Split(",")) {
    $ws.Range($addr).Clear() | Out-Null
}

# Restore the active selection left behind after the refresh/edit.
$ws.Range("C14").Select() | Out-Null
